# Fruta / hortaliza, semanal
# A new weekly price entry is prepended to the data table (row 142),
# pushing the existing rows 142-186 down to rows 143-187.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 142, shifting rows 142:186 down to 143:187.
$ws.Rows.Item(142).Insert()

# Populate the newly inserted row 142 with this week's entry.
$ws.Range("A142").Value = 2
$ws.Range("B142").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C142").Value = "Coquimbo"
$ws.Range("D142").Value = 44924
$ws.Range("E142").Value = 4
$ws.Range("F142").Value = 100112043
$ws.Range("G142").Value = "Pepino ensalada"
$ws.Range("H142").Value = "Sin especificar"
$ws.Range("I142").Value = "Primera"
$ws.Range("J142").Value = 400
$ws.Range("K142").Value = 13000
$ws.Range("L142").Value = 14000
$ws.Range("M142").Value = 13500
$ws.Range("N142").Value = "$/caja 70 unidades"
$ws.Range("O142").Value = "Provincia de Limarí"
$ws.Range("P142").Value = 193
$ws.Range("Q142").Value = 70
$ws.Range("R142").Value = "Hortaliza"
